$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1805.4615
$ws.Range("J19").Value = 327.625
$ws.Range("L19").Value = 327.625
$ws.Range("N19").Value = -677.625
$ws.Range("H28").Value = 599.6
$ws.Range("I28").Value = 199.66667
$ws.Range("J28").Value = 1199.5
$ws.Range("K28").Value = 199.66667
$ws.Range("L28").Value = 1199.5
$ws.Range("M28").Value = 285.33333
$ws.Range("N28").Value = -2169.5
$ws.Range("H40").Value = 1141.9412
$ws.Range("I40").Value = 1015.1
$ws.Range("J40").Value = 1323.1428
$ws.Range("K40").Value = 1015.1
$ws.Range("L40").Value = 1323.1428
$ws.Range("M40").Value = -840.1
$ws.Range("N40").Value = -1673.1428
$ws.Range("H62").Value = 2659.2083
$ws.Range("I62").Value = 2018.5714
$ws.Range("J62").Value = 3556.1
$ws.Range("K62").Value = 2018.5714
$ws.Range("L62").Value = 3556.1
$ws.Range("M62").Value = -1394.5714
$ws.Range("N62").Value = -4804.1
$ws.Range("H64").Value = 3217.96
$ws.Range("I64").Value = 2823.077
$ws.Range("J64").Value = 3645.75
$ws.Range("K64").Value = 2823.077
$ws.Range("L64").Value = 3645.75
$ws.Range("M64").Value = -2575.077
$ws.Range("N64").Value = -4141.75
$ws.Range("H65").Value = 2659.2083
$ws.Range("I65").Value = 2018.5714
$ws.Range("J65").Value = 3556.1
$ws.Range("K65").Value = 10092.857
$ws.Range("L65").Value = 17780.5
$ws.Range("M65").Value = -6972.857
$ws.Range("N65").Value = -24020.5
$ws.Range("H67").Value = 3217.96
$ws.Range("I67").Value = 2823.077
$ws.Range("J67").Value = 3645.75
$ws.Range("K67").Value = 2823.077
$ws.Range("L67").Value = 3645.75
$ws.Range("M67").Value = -1965.077
$ws.Range("N67").Value = -5361.75
$ws.Range("H74").Value = 13893556
$ws.Range("J74").Value = 17861428
$ws.Range("L74").Value = 17861428
$ws.Range("N74").Value = -17863300
$ws.Range("H77").Value = 13893556
$ws.Range("J77").Value = 17861428
$ws.Range("L77").Value = 89307140
$ws.Range("N77").Value = -89316500
$ws.Range("H106").Value = 13335398
$ws.Range("I106").Value = 17545438
$ws.Range("K106").Value = 17545438
$ws.Range("M106").Value = -17544807
$ws.Range("H116").Value = 5083.25
$ws.Range("I116").Value = 2598.8
$ws.Range("J116").Value = 6857.857
$ws.Range("K116").Value = 2598.8
$ws.Range("L116").Value = 6857.857
$ws.Range("M116").Value = 843.1999999999998
$ws.Range("N116").Value = -13741.857
$ws.Range("H132").Value = 3132.3928
$ws.Range("I132").Value = 3315.6538
$ws.Range("K132").Value = 9946.9614
$ws.Range("M132").Value = -7416.9614
$ws.Range("H137").Value = 1833.05
$ws.Range("I137").Value = 1516.3125
$ws.Range("J137").Value = 3100
$ws.Range("K137").Value = 4548.9375
$ws.Range("L137").Value = 9300
$ws.Range("M137").Value = -1998.9375
$ws.Range("N137").Value = -14400
$ws.Range("H141").Value = 2053.7058
$ws.Range("I141").Value = 1853.8667
$ws.Range("J141").Value = 3552.5
$ws.Range("K141").Value = 5561.6001
$ws.Range("L141").Value = 10657.5
$ws.Range("M141").Value = -381.6000999999997
$ws.Range("N141").Value = -21017.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1841.9459
$ws.Range("I45").Value = 2375
$ws.Range("J45").Value = 1478.5
$ws.Range("K45").Value = 2375
$ws.Range("L45").Value = 1478.5
$ws.Range("M45").Value = -1998
$ws.Range("N45").Value = -2232.5
$ws.Range("H74").Value = 55556444
$ws.Range("I74").Value = 66667164
$ws.Range("K74").Value = 66667164
$ws.Range("M74").Value = -66666290
$ws.Range("H77").Value = 55556444
$ws.Range("I77").Value = 66667164
$ws.Range("K77").Value = 333335820
$ws.Range("M77").Value = -333331452
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20698.084
$ws.Range("J82").Value = 49108.75
$ws.Range("L82").Value = 49108.75
$ws.Range("N82").Value = -49874.75
$ws.Range("H85").Value = 20698.084
$ws.Range("J85").Value = 49108.75
$ws.Range("L85").Value = 49108.75
$ws.Range("N85").Value = -51760.75
$ws.Range("H94").Value = 4000
$ws.Range("I94").Value = 4000
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 4000
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -3549
$ws.Range("N94").Value = -4902
$ws.Range("H110").Value = 45466
$ws.Range("J110").Value = 45466
$ws.Range("L110").Value = 45466
$ws.Range("N110").Value = -53646
$ws.Range("H134").Value = 4560.12
$ws.Range("I134").Value = 4963.773
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 14891.319
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -12356.319
$ws.Range("N134").Value = -9870

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 22.7
$ws.Range("I7").Value = 15
$ws.Range("K7").Value = 15
$ws.Range("M7").Value = 98
$ws.Range("H16").Value = 997.8461
$ws.Range("I16").Value = 1020.375
$ws.Range("J16").Value = 961.8
$ws.Range("K16").Value = 1020.375
$ws.Range("L16").Value = 961.8
$ws.Range("M16").Value = -733.375
$ws.Range("N16").Value = -1535.8
$ws.Range("H22").Value = 263.46667
$ws.Range("I22").Value = 86.36364
$ws.Range("J22").Value = 750.5
$ws.Range("K22").Value = 86.36364
$ws.Range("L22").Value = 750.5
$ws.Range("M22").Value = 263.63636
$ws.Range("N22").Value = -1450.5
$ws.Range("H31").Value = 18366.191
$ws.Range("I31").Value = 39376.625
$ws.Range("K31").Value = 39376.625
$ws.Range("M31").Value = -39081.625
$ws.Range("H34").Value = 18366.191
$ws.Range("I34").Value = 39376.625
$ws.Range("K34").Value = 39376.625
$ws.Range("M34").Value = -39174.625
$ws.Range("H62").Value = 7500
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 7500
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H99").Value = 3950
$ws.Range("I99").Value = 3057.8948
$ws.Range("K99").Value = 3057.8948
$ws.Range("M99").Value = -1559.8948
$ws.Range("H113").Value = 997.8461
$ws.Range("I113").Value = 1020.375
$ws.Range("J113").Value = 961.8
$ws.Range("K113").Value = 1020.375
$ws.Range("L113").Value = 961.8
$ws.Range("M113").Value = 1149.625
$ws.Range("N113").Value = -5301.8
$ws.Range("H126").Value = 3950
$ws.Range("I126").Value = 3057.8948
$ws.Range("K126").Value = 9173.6844
$ws.Range("M126").Value = -6703.6844
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2727363.8
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 6000080
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 18000240
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -18000464
$ws.Range("H113").Value = 542.5
$ws.Range("I113").Value = 448.2
$ws.Range("J113").Value = 699.6667
$ws.Range("K113").Value = 1344.6
$ws.Range("L113").Value = 2099.0001
$ws.Range("M113").Value = 825.4000000000001
$ws.Range("N113").Value = -6439.0001
$ws.Range("H131").Value = 724.38
$ws.Range("J131").Value = 724.38
$ws.Range("L131").Value = 2173.14
$ws.Range("N131").Value = -12253.14

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 50003500
$ws.Range("J52").Value = 50003500
$ws.Range("L52").Value = 50003500
$ws.Range("N52").Value = -50004018
$ws.Range("H113").Value = 3269.2307
$ws.Range("I113").Value = 2555.5557
$ws.Range("J113").Value = 4875
$ws.Range("K113").Value = 2555.5557
$ws.Range("L113").Value = 4875
$ws.Range("M113").Value = -385.5556999999999
$ws.Range("N113").Value = -9215

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 800.6842
$ws.Range("I46").Value = 786.46155
$ws.Range("K46").Value = 786.46155
$ws.Range("M46").Value = -598.46155
$ws.Range("H61").Value = 5588.857
$ws.Range("I61").Value = 3161.111
$ws.Range("K61").Value = 3161.111
$ws.Range("M61").Value = -2959.111
$ws.Range("H113").Value = 5588.857
$ws.Range("I113").Value = 3161.111
$ws.Range("K113").Value = 3161.111
$ws.Range("M113").Value = -991.1109999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1575.1666
$ws.Range("I81").Value = 1360
$ws.Range("J81").Value = 2651
$ws.Range("K81").Value = 2720
$ws.Range("L81").Value = 5302
$ws.Range("M81").Value = -1659
$ws.Range("N81").Value = -7424
$ws.Range("H84").Value = 1575.1666
$ws.Range("I84").Value = 1360
$ws.Range("J84").Value = 2651
$ws.Range("K84").Value = 13600
$ws.Range("L84").Value = 26510
$ws.Range("M84").Value = -8296
$ws.Range("N84").Value = -37118
$ws.Range("H113").Value = 1378
$ws.Range("I113").Value = 1485.8
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 4457.4
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = -2287.4
$ws.Range("N113").Value = -5240
$ws.Range("H136").Value = 30304854
$ws.Range("I136").Value = 43479800
$ws.Range("J136").Value = 2482.4
$ws.Range("K136").Value = 130439400
$ws.Range("L136").Value = 7447.200000000001
$ws.Range("M136").Value = -130436850
$ws.Range("N136").Value = -12547.2
